$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$dCell = $ws.Cells.Item(2, 4)
$dCell.Value = "'27.535.76"
$dCell.Style = "Normal"
$ws.Range("E2").Value = "  +4.73%  "

$dCell = $ws.Cells.Item(3, 4)
$dCell.Value = "'1.841.53"
$dCell.Style = "Normal"
$ws.Range("E3").Value = "  +3.94%  "

$dCell = $ws.Cells.Item(4, 4)
$dCell.Value = "'1.027"
$dCell.Style = "Normal"
$ws.Range("E4").Value = "  +2.48%  "

$dCell = $ws.Cells.Item(5, 4)
$dCell.Value = "'319.19"
$dCell.Style = "Normal"
$ws.Range("E5").Value = "  +4.41%  "

$dCell = $ws.Cells.Item(6, 4)
$dCell.Value = "'1.023"
$dCell.Style = "Normal"
$ws.Range("E6").Value = "  +2.10%  "

$dCell = $ws.Cells.Item(7, 4)
$dCell.Value = "'0.4380"
$dCell.Style = "Normal"
$ws.Range("E7").Value = "  +3.57%  "

$dCell = $ws.Cells.Item(8, 4)
$dCell.Value = "'0.3737"
$dCell.Style = "Normal"
$ws.Range("E8").Value = "  +3.73%  "

$dCell = $ws.Cells.Item(9, 4)
$dCell.Value = "'0.07386"
$dCell.Style = "Normal"
$ws.Range("E9").Value = "  +3.71%  "

$dCell = $ws.Cells.Item(10, 4)
$dCell.Value = "'0.8765"
$dCell.Style = "Normal"
$ws.Range("E10").Value = "  +4.78%  "

$dCell = $ws.Cells.Item(11, 4)
$dCell.Value = "'21.52"
$dCell.Style = "Normal"
$ws.Range("E11").Value = "  +5.76%  "

$dCell = $ws.Cells.Item(12, 4)
$dCell.Value = "'1.860.66"
$dCell.Style = "Normal"
$ws.Range("E12").Value = "  +6.04%  "

$dCell = $ws.Cells.Item(13, 4)
$dCell.Value = "'5.489"
$dCell.Style = "Normal"
$ws.Range("E13").Value = "  +4.78%  "

$dCell = $ws.Cells.Item(14, 4)
$dCell.Value = "'6.681"
$dCell.Style = "Normal"
$ws.Range("E14").Value = "  +3.71%  "

$dCell = $ws.Cells.Item(15, 4)
$dCell.Value = "'0.07126"
$dCell.Style = "Normal"
$ws.Range("E15").Value = "  +3.46%  "

$dCell = $ws.Cells.Item(16, 4)
$dCell.Value = "'82.72"
$dCell.Style = "Normal"
$ws.Range("E16").Value = "  +4.92%  "

$dCell = $ws.Cells.Item(17, 4)
$dCell.Value = "'1.030"
$dCell.Style = "Normal"
$ws.Range("E17").Value = "  +2.74%  "

$dCell = $ws.Cells.Item(18, 4)
$dCell.Value = "'0.000009000"
$dCell.Style = "Normal"
$ws.Range("E18").Value = "  +4.27%  "

$ws.Range("E19").Value = "  +2.47%  "

$dCell = $ws.Cells.Item(20, 4)
$dCell.Value = "'15.40"
$dCell.Style = "Normal"
$ws.Range("E20").Value = "  +3.44%  "

$dCell = $ws.Cells.Item(21, 4)
$dCell.Value = "'27.532.37"
$dCell.Style = "Normal"
$ws.Range("E21").Value = "  +4.70%  "

$dCell = $ws.Cells.Item(22, 4)
$dCell.Value = "'5.246"
$dCell.Style = "Normal"
$ws.Range("E22").Value = "  +2.83%  "

$dCell = $ws.Cells.Item(23, 4)
$dCell.Value = "'11.19"
$dCell.Style = "Normal"
$ws.Range("E23").Value = "  +1.85%  "

$dCell = $ws.Cells.Item(24, 4)
$dCell.Value = "'2.077.75"
$dCell.Style = "Normal"
$ws.Range("E24").Value = "  +5.19%  "

$dCell = $ws.Cells.Item(25, 4)
$dCell.Value = "'157.12"
$dCell.Style = "Normal"
$ws.Range("E25").Value = "  +3.43%  "

$dCell = $ws.Cells.Item(26, 4)
$dCell.Value = "'1.922"
$dCell.Style = "Normal"
$ws.Range("E26").Value = "  +7.10%  "

$dCell = $ws.Cells.Item(27, 4)
$dCell.Value = "'18.72"
$dCell.Style = "Normal"
$ws.Range("E27").Value = "  +3.97%  "

$dCell = $ws.Cells.Item(28, 4)
$dCell.Value = "'5.261"
$dCell.Style = "Normal"
$ws.Range("E28").Value = "  +4.06%  "

$dCell = $ws.Cells.Item(29, 4)
$dCell.Value = "'1.941"
$dCell.Style = "Normal"
$ws.Range("E29").Value = "  +5.87%  "

$dCell = $ws.Cells.Item(30, 4)
$dCell.Value = "'116.39"
$dCell.Style = "Normal"
$ws.Range("E30").Value = "  +2.02%  "

$dCell = $ws.Cells.Item(31, 4)
$dCell.Value = "'0.09086"
$dCell.Style = "Normal"
$ws.Range("E31").Value = "  +3.07%  "

$dCell = $ws.Cells.Item(32, 4)
$dCell.Value = "'1.208"
$dCell.Style = "Normal"
$ws.Range("E32").Value = "  +8.34%  "

$dCell = $ws.Cells.Item(33, 4)
$dCell.Value = "'0.7674"
$dCell.Style = "Normal"
$ws.Range("E33").Value = "  +5.91%  "

$dCell = $ws.Cells.Item(34, 4)
$dCell.Value = "'4.499"
$dCell.Style = "Normal"
$ws.Range("E34").Value = "  +4.36%  "

$dCell = $ws.Cells.Item(35, 4)
$dCell.Value = "'2.869"
$dCell.Style = "Normal"
$ws.Range("E35").Value = "  +4.92%  "

$dCell = $ws.Cells.Item(36, 4)
$dCell.Value = "'1.027"
$dCell.Style = "Normal"
$ws.Range("E36").Value = "  +2.61%  "

$dCell = $ws.Cells.Item(37, 4)
$dCell.Value = "'1.145"
$dCell.Style = "Normal"
$ws.Range("E37").Value = "  +3.72%  "

$dCell = $ws.Cells.Item(38, 4)
$dCell.Value = "'0.01973"
$dCell.Style = "Normal"
$ws.Range("E38").Value = "  +4.72%  "

$dCell = $ws.Cells.Item(39, 4)
$dCell.Value = "'0.05252"
$dCell.Style = "Normal"
$ws.Range("E39").Value = "  +2.92%  "

$dCell = $ws.Cells.Item(40, 4)
$dCell.Value = "'0.5180"
$dCell.Style = "Normal"
$ws.Range("E40").Value = "  +5.41%  "

$dCell = $ws.Cells.Item(41, 4)
$dCell.Value = "'2.783"
$dCell.Style = "Normal"
$ws.Range("E41").Value = "  +7.35%  "

$dCell = $ws.Cells.Item(42, 4)
$dCell.Value = "'0.1668"
$dCell.Style = "Normal"
$ws.Range("E42").Value = "  +3.66%  "

$dCell = $ws.Cells.Item(43, 4)
$dCell.Value = "'6.631"
$dCell.Style = "Normal"

$dCell = $ws.Cells.Item(44, 4)
$dCell.Value = "'8.551"
$dCell.Style = "Normal"
$ws.Range("E44").Value = "  +6.57%  "

$dCell = $ws.Cells.Item(45, 4)
$dCell.Value = "'108.96"
$dCell.Style = "Normal"
$ws.Range("E45").Value = "  +4.08%  "

$dCell = $ws.Cells.Item(46, 4)
$dCell.Value = "'10.54"
$dCell.Style = "Normal"
$ws.Range("E46").Value = "  +3.82%  "

$dCell = $ws.Cells.Item(47, 4)
$dCell.Value = "'1.027"
$dCell.Style = "Normal"
$ws.Range("E47").Value = "  +2.57%  "

$dCell = $ws.Cells.Item(48, 4)
$dCell.Value = "'1.706"
$dCell.Style = "Normal"
$ws.Range("E48").Value = "  +5.56%  "

$dCell = $ws.Cells.Item(49, 4)
$dCell.Value = "'0.4653"
$dCell.Style = "Normal"
$ws.Range("E49").Value = "  +4.81%  "

$dCell = $ws.Cells.Item(50, 4)
$dCell.Value = "'1.898"
$dCell.Style = "Normal"
$ws.Range("E50").Value = "  +12.19%  "

$dCell = $ws.Cells.Item(51, 4)
$dCell.Value = "'0.06340"
$dCell.Style = "Normal"
$ws.Range("E51").Value = "  +2.80%  "
